$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.926.69'
$ws.Range("E2").Value = '  +1.72%  '

$ws.Range("D3").Value = '1.645.39'
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("E4").Value = '  -0.12%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '213.51'
$c.Style = "Normal"

$ws.Range("E6").Value = '  +0.23%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.40'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.94%  '

$ws.Range("E9").Value = '  +1.61%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0614'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("E11").Value = '  -1.62%  '

$ws.Range("D12").Value = '1.878.82'
$ws.Range("E12").Value = '  +1.87%  '

$ws.Range("D13").Value = '1.643.71'
$ws.Range("E13").Value = '  +1.75%  '

$ws.Range("E14").Value = '  +1.10%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.55'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("D17").Value = '27.936.90'

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '231.26'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.72%  '

$ws.Range("E19").Value = '  +1.14%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.00%  '

$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("E22").Value = '  +4.79%  '

$ws.Range("E23").Value = '  +2.07%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +3.37%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '152.36'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.55%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.92'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.95%  '

$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("E28").Value = '  +1.48%  '

$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("E30").Value = '  +1.63%  '

$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("E32").Value = '  +2.21%  '

$ws.Range("D33").Value = '1.442.06'
$ws.Range("E33").Value = '  -1.85%  '

$ws.Range("E35").Value = '  +1.79%  '

$ws.Range("E36").Value = '  -0.20%  '

$ws.Range("E37").Value = '  +3.42%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.935'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.22%  '

$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("E40").Value = '  +0.61%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '69.12'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.49%  '

$ws.Range("E42").Value = '  +3.46%  '

$ws.Range("E43").Value = '  -0.12%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.45'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.29%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.82'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +5.74%  '

$ws.Range("E46").Value = '  +3.55%  '

$ws.Range("E47").Value = '  +0.71%  '

$ws.Range("D48").Value = '1.787.23'
$ws.Range("E48").Value = '  +1.57%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '89.03'
$c.Style = "Normal"

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.101'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0508'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.15%  '
